$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append, continuing from the existing last row (229)
$data = @(
    @(44304, 0, 2, 57.75339301183945),
    @(44305, 0, 2, 57.75339301183945),
    @(44306, 0, 2, 57.75339301183945),
    @(44307, 1, 3, 86.63008951775916)
)

$lastRow = 229
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy formatting of the last existing row's date cell down into the new
    # row so the new cell picks up the same style (the date format on column A).
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
